# Apply crypto price/volume updates from the Jan 18 2024 GitHub Actions refresh.
# Source cells are text-formatted (t="inlineStr"/shared string), so force
# NumberFormat to "@" (Text) before assigning to preserve exact string
# representations (e.g. trailing zeros like "1.00", "0.110", "155.60").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.748.63"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.530.90"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.99%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.88"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.01"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.82%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.16%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.517"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.34"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0802"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.29%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.110"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.24"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.917.51"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.26%  "
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.26"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -6.23%  "
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.515.66"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.813"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.58%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.698.29"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.60"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -3.08%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0942"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.09"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "242.38"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.26%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.00"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.49%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.55"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -3.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.26"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -4.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.02"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.69"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -5.94%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +4.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "155.60"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.72"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0784"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.46%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.69%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.14"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.49%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -4.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "17.66"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.70%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.60%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.87%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.24"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.28%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "22.36"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.58%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.26%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.031.68"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.49%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0297"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.23"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.88"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.21%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.769.66"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.29%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "80.45"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.60%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.56%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.16"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.83%  "
